$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (D2): part rating corrected from 80V (8TQ080) to 100V (8TQ100)
$ws.Range("A7").Value = "'8TQ100"
$ws.Range("B7").Value = "'8A, 100V, TO-220AC, Schottky Rect"
$ws.Range("D7").Value = "'8TQ100"

# Row 25: heatsink entry (V5629G / TO-220F Heatsink) removed - clear the row but keep formatting
$ws.Range("A25:F25").ClearContents()

# Row 26: renumber label (was HeatSink2, now HeatSink1); part/price untouched
$ws.Range("A26").Value = "HeatSink1"

# Row 27: renumber label (was HeatSink3, now HeatSink2); quantities/price combined (doubled)
$ws.Range("A27").Value = "HeatSink2"
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 0.56000000000000005

# Update the active selection to match the saved view
$null = $ws.Range("B30").Select()
